$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.130.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.783.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'629.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'164.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.49%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.781.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.30%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.28%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.09%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'6.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.70%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'35.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.14%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.420.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.26%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.796.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'69.173.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'17.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.42%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.30%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'468.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.60%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.86%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'83.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'12.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.10%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.84%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.05%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.931.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.35%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.98%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.06%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.74%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'28.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.29%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.20%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.733.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +7.91%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.28%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E45").Value = "'  -0.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'153.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.94%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +2.87%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'46.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'42.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.61%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.81%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.69%  "
$ws.Range("E51").Style = "Normal"
